$wb = $excel.ActiveWorkbook
$status = $wb.Worksheets.Item("Status")
$songList = $wb.Worksheets.Item("SongList")

# --- Add new rows of pseudo-code / notes to the Status sheet ---

$status.Range("C35").Value = "Player 1 Right"
$status.Range("D35").Value = "if ((playerRound % 2 === 0) && (currentSong === playerSelection) && (playerRound < 9)) {"
$status.Range("D36").Value = '    correctWrong.innerHTML = "Correct Answer!";'
$status.Range("D37").Value = '    playerTurn.innerHTML = "Player 1 Turn";'
$status.Range("D38").Value = "    currentPlayer.score += 1;"
$status.Range("D39").Value = "    playerRound = playerRound + 1;"
$status.Range("D40").Value = "    randomAlbumArray();"
$status.Range("D41").Value = "    newRandomSong();"

$status.Range("C42").Value = "Player 1 Wrong"
$status.Range("D42").Value = "} else if ((playerRound % 2 === 0) && (currentSong !== playerSelection) && (playerRound < 9)) {"
$status.Range("D43").Value = '    correctWrong.innerHTML = "Wrong Answer!";'
$status.Range("D44").Value = '    playerTurn.innerHTML = "Player 1 Turn";'
$status.Range("D45").Value = "    currentPlayer.score += 1;"
$status.Range("D46").Value = "    playerRound = playerRound + 1;"
$status.Range("D47").Value = "    randomAlbumArray();"
$status.Range("D48").Value = "    newRandomSong();"

$status.Range("C49").Value = "Player 2 Right "
$status.Range("D49").Value = "} else if ((currentSong === playerSelection) && (playerRound < 9)) {"
$status.Range("D50").Value = '    correctWrong.innerHTML = "Correct Answer!";'
$status.Range("D51").Value = '    playerTurn.innerHTML = "Player 2 Turn";'
$status.Range("D52").Value = "    currentPlayer.score += 1;"
$status.Range("D53").Value = "    playerRound = playerRound + 1;"
$status.Range("D54").Value = "    randomAlbumArray();"
$status.Range("D55").Value = "    newRandomSong();"

$status.Range("C56").Value = "Player 2 Wrong"
$status.Range("D56").Value = "} else if ((currentSong !== playerSelection) && (playerRound < 9)) {"
$status.Range("D57").Value = '    correctWrong.innerHTML = "Wrong Answer!";'
$status.Range("D58").Value = '    playerTurn.innerHTML = "Player 2 Turn";'
$status.Range("D59").Value = "    currentPlayer.score += 1;"
$status.Range("D60").Value = "    playerRound = playerRound + 1;"
$status.Range("D61").Value = "    randomAlbumArray();"
$status.Range("D62").Value = "    newRandomSong();"

$status.Range("C63").Value = "Game End- Player 1 Wins"
$status.Range("D63").Value = "} else if (playerRound =9){"
$status.Range("D64").Value = 'playerTurn.innerHTML = "Player 1 Wins!"'
$status.Range("D65").Value = "}"

$status.Range("C66").Value = "Game End- Player 2 Wins"
$status.Range("D66").Value = "} else {"
$status.Range("D67").Value = 'playerTurn.innerHTML = "Player 2 Wins!"'
$status.Range("D68").Value = "}"

$status.Range("B73").Value = "need to start by listing Player 1 Turn on load"
$status.Range("B75").Value = "when game finishes list:"
$status.Range("B77").Value = "Game Over "
$status.Range("B78").Value = "Player 1 Wins 5 to 3"
$status.Range("B79").Value = "Player 2 Wins 4 to 2"

$status.Columns("C").ColumnWidth = 24.8

# --- View / selection changes ---

# SongList becomes the active/selected sheet, scrolled down, with a new selection
$songList.Activate()
$win = $wb.Windows.Item(1)
$win.FreezePanes = $false
$songList.Range("A2").Select()
$win.FreezePanes = $true
$songList.Range("E37").Select()

# Status sheet keeps a new (non-active) selection
$status.Range("B70").Select()

$songList.Activate()
